$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B21").Value = "Spacing"
$ws.Range("C21").Value = 0.41699999999999998
$ws.Range("D21").Value = 0.46899999999999997
$ws.Range("E21").Value = 0.59099999999999997
$ws.Range("F21").Value = 0.46899999999999997
$ws.Range("G21").Value = 0.41699999999999998

$ws.Range("B22").Value = "Phase"
$ws.Range("C22").Value = 34.799999999999997
$ws.Range("D22").Value = 239.5
$ws.Range("E22").Value = 112.8
$ws.Range("F22").Value = 313.7
$ws.Range("G22").Value = 218.7
$ws.Range("H22").Value = 106.7

$ws.Range("B23").Value = "Amplitude"
$ws.Range("C23").Value = 0.73399999999999999
$ws.Range("D23").Value = 0.98099999999999998
$ws.Range("E23").Value = 0.98799999999999999
$ws.Range("F23").Value = 0.59499999999999997
$ws.Range("G23").Value = 0.93799999999999994
$ws.Range("H23").Value = 1

$ws.Range("B24").Value = "Fitness"
$ws.Range("C24").Value = -4.0476846462134679
